$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 0.225969755229471
$ws.Cells.Item(2, 4).Value = 0.006298024302926919
$ws.Cells.Item(2, 5).Value = 0.1393848103014683
$ws.Cells.Item(2, 6).Value = 0.7533793398549165
$ws.Cells.Item(2, 7).Value = 0.6201500300061298
$ws.Cells.Item(2, 8).Value = 0.6454766957521656
$ws.Cells.Item(2, 9).Value = 0.7538371318557182
$ws.Cells.Item(2, 12).Value = 0.1684885293036018
$ws.Cells.Item(2, 14).Value = 2.873581687696344
$ws.Cells.Item(2, 15).Value = 2.513635408739276

$ws.Cells.Item(3, 3).Value = 0.2233335232342313
$ws.Cells.Item(3, 4).Value = 0.006197097785561212
$ws.Cells.Item(3, 5).Value = 0.1363420673065541
$ws.Cells.Item(3, 6).Value = 0.7228535965416398
$ws.Cells.Item(3, 7).Value = 0.5877630113456149
$ws.Cells.Item(3, 8).Value = 0.6345060536392282
$ws.Cells.Item(3, 9).Value = 0.7294046322488512
$ws.Cells.Item(3, 12).Value = 0.1634208008851701
$ws.Cells.Item(3, 14).Value = 2.562605684679454
$ws.Cells.Item(3, 15).Value = 2.42162939103568

$ws.Cells.Item(4, 3).Value = 0.221825637815698
$ws.Cells.Item(4, 4).Value = 0.006133699267213188
$ws.Cells.Item(4, 5).Value = 0.1345519741412922
$ws.Cells.Item(4, 6).Value = 0.7045593754398993
$ws.Cells.Item(4, 7).Value = 0.5682508425061883
$ws.Cells.Item(4, 8).Value = 0.6281106818067173
$ws.Cells.Item(4, 9).Value = 0.7148438111719386
$ws.Cells.Item(4, 12).Value = 0.1604095489491826
$ws.Cells.Item(4, 14).Value = 2.371325805375761
$ws.Cells.Item(4, 15).Value = 2.366627523460778

$ws.Cells.Item(5, 3).Value = 0.2212390354097806
$ws.Cells.Item(5, 4).Value = 0.006107505335972263
$ws.Cells.Item(5, 5).Value = 0.1338421385255373
$ws.Cells.Item(5, 6).Value = 0.6972170454626365
$ws.Cells.Item(5, 7).Value = 0.5603932113413777
$ws.Cells.Item(5, 8).Value = 0.6255902006696346
$ws.Cells.Item(5, 9).Value = 0.7090209246642019
$ws.Cells.Item(5, 12).Value = 0.1592076312726221
$ws.Cells.Item(5, 14).Value = 2.293303068607429
$ws.Cells.Item(5, 15).Value = 2.344588144831079

$ws.Cells.Item(6, 3).Value = 0.2211433145573523
$ws.Cells.Item(6, 4).Value = 0.006103134226124851
$ws.Cells.Item(6, 5).Value = 0.1337254568156716
$ws.Cells.Item(6, 6).Value = 0.6960046614146194
$ws.Cells.Item(6, 7).Value = 0.5590941131201674
$ws.Cells.Item(6, 8).Value = 0.6251768526705348
$ws.Cells.Item(6, 9).Value = 0.7080607263994665
$ws.Cells.Item(6, 12).Value = 0.1590095740485324
$ws.Cells.Item(6, 14).Value = 2.280343261403573
$ws.Cells.Item(6, 15).Value = 2.340951120064545

$ws.Cells.Item(7, 3).Value = 0.2218176138079286
$ws.Cells.Item(7, 4).Value = 0.006133347457309313
$ws.Cells.Item(7, 5).Value = 0.1345423215372392
$ws.Cells.Item(7, 6).Value = 0.7044598979951644
$ws.Cells.Item(7, 7).Value = 0.5681444924857146
$ws.Cells.Item(7, 8).Value = 0.6280763427505178
$ws.Cells.Item(7, 9).Value = 0.7147648333137866
$ws.Cells.Item(7, 12).Value = 0.1603932375217596
$ws.Cells.Item(7, 14).Value = 2.370273851395496
$ws.Cells.Item(7, 15).Value = 2.366328777579895

$ws.Cells.Item(8, 3).Value = 0.2250378015565531
$ws.Cells.Item(8, 4).Value = 0.006263521966417329
$ws.Cells.Item(8, 5).Value = 0.1383194406000712
$ws.Cells.Item(8, 6).Value = 0.7427608516765218
$ws.Cells.Item(8, 7).Value = 0.6089053112631149
$ws.Cells.Item(8, 8).Value = 0.6416233275604952
$ws.Cells.Item(8, 9).Value = 0.7453212505271836
$ws.Cells.Item(8, 12).Value = 0.1667203203116685
$ws.Cells.Item(8, 14).Value = 2.766433886209654
$ws.Cells.Item(8, 15).Value = 2.481602170323981

$ws.Cells.Item(9, 3).Value = 0.2322312903143882
$ws.Cells.Item(9, 4).Value = 0.006507425466280026
$ws.Cells.Item(9, 5).Value = 0.1463477039678622
$ws.Cells.Item(9, 6).Value = 0.821439386749816
$ws.Cells.Item(9, 7).Value = 0.6918155861135062
$ws.Cells.Item(9, 8).Value = 0.6708929033100617
$ws.Cells.Item(9, 9).Value = 0.8087480417135993
$ws.Cells.Item(9, 12).Value = 0.1799267614938174
$ws.Cells.Item(9, 14).Value = 3.540180268007646
$ws.Cells.Item(9, 15).Value = 2.719509808838893

$ws.Cells.Item(10, 3).Value = 0.2380526478987974
$ws.Cells.Item(10, 4).Value = 0.006679668274898631
$ws.Cells.Item(10, 5).Value = 0.1526272775549558
$ws.Cells.Item(10, 6).Value = 0.8814422773569959
$ws.Cells.Item(10, 7).Value = 0.7545720582328102
$ws.Cells.Item(10, 8).Value = 0.6940505126206915
$ws.Cells.Item(10, 9).Value = 0.8575018524364708
$ws.Cells.Item(10, 12).Value = 0.1901217526643251
$ws.Cells.Item(10, 14).Value = 4.10623028343673
$ws.Cells.Item(10, 15).Value = 2.901595926180164

$ws.Cells.Item(11, 3).Value = 0.240817580780444
$ws.Cells.Item(11, 4).Value = 0.006756512684612304
$ws.Cells.Item(11, 5).Value = 0.1555673673806126
$ws.Cells.Item(11, 6).Value = 0.9092213429262728
$ws.Cells.Item(11, 7).Value = 0.7835277731552139
$ws.Cells.Item(11, 8).Value = 0.7049457602524853
$ws.Cells.Item(11, 9).Value = 0.8801529968835524
$ws.Cells.Item(11, 12).Value = 0.1948677723271857
$ws.Cells.Item(11, 14).Value = 4.363110593465422
$ws.Cells.Item(11, 15).Value = 2.986031532157085

$ws.Cells.Item(12, 3).Value = 0.2418813775621942
$ws.Cells.Item(12, 4).Value = 0.006785393981141752
$ws.Cells.Item(12, 5).Value = 0.15669273642186
$ws.Cells.Item(12, 6).Value = 0.9198103212405329
$ws.Cells.Item(12, 7).Value = 0.7945514864979657
$ws.Cells.Item(12, 8).Value = 0.709123415987591
$ws.Cells.Item(12, 9).Value = 0.888798581005787
$ws.Cells.Item(12, 12).Value = 0.1966806036692077
$ws.Cells.Item(12, 14).Value = 4.460285735714251
$ws.Cells.Item(12, 15).Value = 3.018236523349117

$ws.Cells.Item(13, 3).Value = 0.2416515239707735
$ws.Cells.Item(13, 4).Value = 0.006779183587617865
$ws.Cells.Item(13, 5).Value = 0.156449833170214
$ws.Cells.Item(13, 6).Value = 0.9175266955689949
$ws.Cells.Item(13, 7).Value = 0.7921747134234352
$ws.Cells.Item(13, 8).Value = 0.7082213766727534
$ws.Cells.Item(13, 9).Value = 0.8869335699643841
$ws.Cells.Item(13, 12).Value = 0.1962894826810668
$ws.Cells.Item(13, 14).Value = 4.439361943450422
$ws.Cells.Item(13, 15).Value = 3.011290318680267

$ws.Cells.Item(14, 3).Value = 0.2409047638335551
$ws.Cells.Item(14, 4).Value = 0.006758893140659694
$ws.Cells.Item(14, 5).Value = 0.1556597111588545
$ws.Cells.Item(14, 6).Value = 0.9100911079303842
$ws.Cells.Item(14, 7).Value = 0.7844335207503832
$ws.Cells.Item(14, 8).Value = 0.7052884191755311
$ws.Cells.Item(14, 9).Value = 0.880862909143076
$ws.Cells.Item(14, 12).Value = 0.195016601795146
$ws.Cells.Item(14, 14).Value = 4.371107314139522
$ws.Cells.Item(14, 15).Value = 2.988676422814763

$ws.Cells.Item(15, 3).Value = 0.2404495360859329
$ws.Cells.Item(15, 4).Value = 0.006746436237797582
$ws.Cells.Item(15, 5).Value = 0.1551773045706852
$ws.Cells.Item(15, 6).Value = 0.9055456700210982
$ws.Cells.Item(15, 7).Value = 0.7796994822956833
$ws.Cells.Item(15, 8).Value = 0.7034986516532911
$ws.Cells.Item(15, 9).Value = 0.8771533240768434
$ws.Cells.Item(15, 12).Value = 0.194238960682469
$ws.Cells.Item(15, 14).Value = 4.329286057409945
$ws.Cells.Item(15, 15).Value = 2.974854862739619

$ws.Cells.Item(16, 3).Value = 0.2378743014693327
$ws.Cells.Item(16, 4).Value = 0.00667461585382334
$ws.Cells.Item(16, 5).Value = 0.1524368170122656
$ws.Cells.Item(16, 6).Value = 0.8796365894037876
$ws.Cells.Item(16, 7).Value = 0.75268795228871
$ws.Cells.Item(16, 8).Value = 0.6933457404824992
$ws.Cells.Item(16, 9).Value = 0.8560310692201085
$ws.Cells.Item(16, 12).Value = 0.1898137718292361
$ws.Cells.Item(16, 14).Value = 4.089429168003562
$ws.Cells.Item(16, 15).Value = 2.896110169809958

$ws.Cells.Item(17, 3).Value = 0.2363243752458573
$ws.Cells.Item(17, 4).Value = 0.006630169039141265
$ws.Cells.Item(17, 5).Value = 0.1507770096440026
$ws.Cells.Item(17, 6).Value = 0.8638661227901281
$ws.Cells.Item(17, 7).Value = 0.7362217736790626
$ws.Cells.Item(17, 8).Value = 0.6872096428431007
$ws.Cells.Item(17, 9).Value = 0.8431944144987256
$ws.Cells.Item(17, 12).Value = 0.1871268257064287
$ws.Cells.Item(17, 14).Value = 3.94211849064385
$ws.Cells.Item(17, 15).Value = 2.848213919115778

$ws.Cells.Item(18, 3).Value = 0.2354438902955138
$ws.Cells.Item(18, 4).Value = 0.006604462358041019
$ws.Cells.Item(18, 5).Value = 0.1498301889176332
$ws.Cells.Item(18, 6).Value = 0.8548408643002148
$ws.Cells.Item(18, 7).Value = 0.7267892042889059
$ws.Cells.Item(18, 8).Value = 0.6837142764616715
$ws.Cells.Item(18, 9).Value = 0.8358556055237187
$ws.Cells.Item(18, 12).Value = 0.1855915530932606
$ws.Cells.Item(18, 14).Value = 3.857331695637754
$ws.Cells.Item(18, 15).Value = 2.820816185126262

$ws.Cells.Item(19, 3).Value = 0.2351476614824861
$ws.Cells.Item(19, 4).Value = 0.006595734149865251
$ws.Cells.Item(19, 5).Value = 0.149510960561912
$ws.Cells.Item(19, 6).Value = 0.8517928770304479
$ws.Cells.Item(19, 7).Value = 0.7236020768632443
$ws.Cells.Item(19, 8).Value = 0.6825366375811655
$ws.Cells.Item(19, 9).Value = 0.8333784499686203
$ws.Cells.Item(19, 12).Value = 0.1850734844567654
$ws.Cells.Item(19, 14).Value = 3.828614786363971
$ws.Cells.Item(19, 15).Value = 2.811565690694863

$ws.Cells.Item(20, 3).Value = 0.2364882300521032
$ws.Cells.Item(20, 4).Value = 0.006634915186365831
$ws.Cells.Item(20, 5).Value = 0.1509528857280316
$ws.Cells.Item(20, 6).Value = 0.8655402061964423
$ws.Cells.Item(20, 7).Value = 0.7379706567757012
$ws.Cells.Item(20, 8).Value = 0.6878593268105533
$ws.Cells.Item(20, 9).Value = 0.8445562924831762
$ws.Cells.Item(20, 12).Value = 0.1874118009699401
$ws.Cells.Item(20, 14).Value = 3.957806003281064
$ws.Cells.Item(20, 15).Value = 2.853296933605463

$ws.Cells.Item(21, 3).Value = 0.2411236501862817
$ws.Cells.Item(21, 4).Value = 0.006764858860631051
$ws.Cells.Item(21, 5).Value = 0.155891462772594
$ws.Cells.Item(21, 6).Value = 0.9122732289703777
$ws.Cells.Item(21, 7).Value = 0.7867056987996648
$ws.Cells.Item(21, 8).Value = 0.7061484924561796
$ws.Cells.Item(21, 9).Value = 0.8826441607062492
$ws.Cells.Item(21, 12).Value = 0.1953900535503124
$ws.Cells.Item(21, 14).Value = 4.391158149571311
$ws.Cells.Item(21, 15).Value = 2.995312397728412

$ws.Cells.Item(22, 3).Value = 0.2442509461545654
$ws.Cells.Item(22, 4).Value = 0.006848513405595114
$ws.Cells.Item(22, 5).Value = 0.1591891827105982
$ws.Cells.Item(22, 6).Value = 0.9432219894342353
$ws.Cells.Item(22, 7).Value = 0.8188997715168114
$ws.Cells.Item(22, 8).Value = 0.718403805864483
$ws.Cells.Item(22, 9).Value = 0.9079336965696427
$ws.Cells.Item(22, 12).Value = 0.2006953513016612
$ws.Cells.Item(22, 14).Value = 4.673791817957863
$ws.Cells.Item(22, 15).Value = 3.089474715421034

$ws.Cells.Item(23, 3).Value = 0.242572907183316
$ws.Cells.Item(23, 4).Value = 0.006803982008188214
$ws.Cells.Item(23, 5).Value = 0.1574227102323036
$ws.Cells.Item(23, 6).Value = 0.9266668589498721
$ws.Cells.Item(23, 7).Value = 0.8016857412671072
$ws.Cells.Item(23, 8).Value = 0.7118352587157233
$ws.Cells.Item(23, 9).Value = 0.8943998444289605
$ws.Cells.Item(23, 12).Value = 0.1978554671861161
$ws.Cells.Item(23, 14).Value = 4.523002190005457
$ws.Cells.Item(23, 15).Value = 3.039095115121142

$ws.Cells.Item(24, 3).Value = 0.2364141183124815
$ws.Cells.Item(24, 4).Value = 0.006632769930758897
$ws.Cells.Item(24, 5).Value = 0.1508733490257868
$ws.Cells.Item(24, 6).Value = 0.8647832241589697
$ws.Cells.Item(24, 7).Value = 0.737179880661273
$ws.Cells.Item(24, 8).Value = 0.6875655038883508
$ws.Cells.Item(24, 9).Value = 0.8439404592562028
$ws.Cells.Item(24, 12).Value = 0.1872829341031803
$ws.Cells.Item(24, 14).Value = 3.950713976768498
$ws.Cells.Item(24, 15).Value = 2.850998471484559

$ws.Cells.Item(25, 3).Value = 0.230191126109716
$ws.Cells.Item(25, 4).Value = 0.006442662608707295
$ws.Cells.Item(25, 5).Value = 0.1441090771026978
$ws.Cells.Item(25, 6).Value = 0.799770333363611
$ws.Cells.Item(25, 7).Value = 0.6690645539529783
$ws.Cells.Item(25, 8).Value = 0.6626847321016385
$ws.Cells.Item(25, 9).Value = 0.7912123505931277
$ws.Cells.Item(25, 12).Value = 0.1762679991569343
$ws.Cells.Item(25, 14).Value = 3.331249627311138
$ws.Cells.Item(25, 15).Value = 2.653873203053706
